$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.21"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'24.01"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'5.351"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05817"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'3.367"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'6.458"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.8092"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.9170"
$ws.Range("D9").Style = "Normal"

$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01070"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1404"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07370"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03173"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03071"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09365"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.846"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001547"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04698"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"

$ws.Range("D19").Value = "'0.006116"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").Value = "'0.001242"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'0.004681"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'0.00008799"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("D23").Value = "'3.592"
$ws.Range("D23").Style = "Normal"

$ws.Range("D28").Value = "'0.0002349"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "27UpBotsUBXTWorstin24h"

$ws.Range("D40").Value = "'0.03841"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006399"
$ws.Range("D41").Style = "Normal"

$ws.Range("D43").Value = "'0.003200"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.009025"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005253"
$ws.Range("D45").Style = "Normal"

$ws.Range("D47").Value = "'0.6854"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'0.001826"
$ws.Range("D48").Style = "Normal"
